# NYPD CompStat weekly report refresh: new crime data collected.
# Updates the report header (issue number / week-covering dates), a batch of
# weekly/28-day/YTD/2-year crime-count + %-change figures, a couple of cells
# that flip from a literal 0 count to the report's "0"/"***.*" placeholder
# text, and the best-fit width of column E after its content changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header: "Volume 31   Number  18" -> "...19" and the week-covering dates.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/6/2024  Through  5/12/2024"

# ---------------------------------------------------------------------------
# Helper: convert a cell that currently holds a plain number into the
# report's literal placeholder text ("0" or "***.*") while keeping the
# surrounding "text" cell style (style 14, as used by neighboring
# placeholder cells such as C15/M15) instead of the numeric style Excel
# would normally infer for a digit-only string.
# ---------------------------------------------------------------------------
function Set-PlaceholderText($addr, $text) {
    $ws.Range("C15").Copy()
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($addr).Formula = '="' + $text + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}

Set-PlaceholderText "D16" "0"
Set-PlaceholderText "E16" "***.*"
Set-PlaceholderText "D18" "0"
Set-PlaceholderText "E18" "***.*"
Set-PlaceholderText "C23" "0"
Set-PlaceholderText "D23" "0"
Set-PlaceholderText "E23" "***.*"

# ---------------------------------------------------------------------------
# Updated weekly/28-day/YTD/2-year counts and %-change figures.
# ---------------------------------------------------------------------------
$ws.Range("L15").Value = -25

$ws.Range("C16").Value = 1
$ws.Range("I16").Value = 31
$ws.Range("K16").Value = 24
$ws.Range("L16").Value = 29.166666666666

$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 66.666666666666
$ws.Range("I17").Value = 80
$ws.Range("J17").Value = 83
$ws.Range("K17").Value = -3.614457831325
$ws.Range("L17").Value = -13.978494623655

$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 22
$ws.Range("K18").Value = 29.411764705882
$ws.Range("L18").Value = -38.888888888888

$ws.Range("C19").Value = 6
$ws.Range("E19").Value = -45.454545454545
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = -21.052631578947
$ws.Range("I19").Value = 129
$ws.Range("J19").Value = 128
$ws.Range("K19").Value = 0.78125
$ws.Range("L19").Value = -27.118644067796

$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -71.428571428571
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -70.588235294117
$ws.Range("I20").Value = 18
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = -55
$ws.Range("L20").Value = -35.714285714285

$ws.Range("C21").Value = 16
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 71
$ws.Range("G21").Value = 77
$ws.Range("H21").Value = -7.792207792207
$ws.Range("I21").Value = 283
$ws.Range("J21").Value = 297
$ws.Range("K21").Value = -4.713804713804
$ws.Range("L21").Value = -21.823204419889

$ws.Range("L23").Value = -72.727272727272

$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = 26.923076923076
$ws.Range("F24").Value = 91
$ws.Range("G24").Value = 125
$ws.Range("H24").Value = -27.2
$ws.Range("I24").Value = 542
$ws.Range("J24").Value = 455
$ws.Range("K24").Value = 19.120879120879
$ws.Range("L24").Value = 22.902494331065

$ws.Range("C25").Value = 21
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 75
$ws.Range("F25").Value = 65
$ws.Range("G25").Value = 81
$ws.Range("H25").Value = -19.753086419753
$ws.Range("I25").Value = 375
$ws.Range("J25").Value = 287
$ws.Range("K25").Value = 30.662020905923
$ws.Range("L25").Value = 60.944206008583

$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -15.384615384615
$ws.Range("F26").Value = 43
$ws.Range("G26").Value = 48
$ws.Range("H26").Value = -10.416666666666
$ws.Range("I26").Value = 223
$ws.Range("J26").Value = 199
$ws.Range("K26").Value = 12.060301507537
$ws.Range("L26").Value = 28.901734104046

$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 6
$ws.Range("K27").Value = -25
$ws.Range("L27").Value = 20

$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -28.571428571428
$ws.Range("I28").Value = 20
$ws.Range("J28").Value = 28
$ws.Range("K28").Value = -28.571428571428
$ws.Range("L28").Value = 33.333333333333

$ws.Range("L31").Value = -40

# ---------------------------------------------------------------------------
# Column E widened (best-fit) now that its contents changed.
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 6.714285714285714
